$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-28
# from 2023-10-05 (serial 45204) to 2023-10-08 (serial 45207)
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
